# Edit script: insert a new date block (44524) of 3 rows before the
# existing row 327, shifting the old rows 327-430 down to 330-433.
# Then populate the 3 newly inserted rows with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at the position currently occupied by row 327.
# This pushes the existing rows 327:430 down to 330:433 and extends
# the used range / dimension accordingly.
$ws.Range("A327:T329").Insert()

# Common (constant-across-sheet) column values reused for the new rows.
$mercadoId = 2
$mercado = "Comercializadora del Agro de Limarí"
$region = "Coquimbo"
$codreg = 4
$tipo = "Fruta"
$productoId = 100102
$producto = "Cítricos"
$categoriaId = 100102003
$categoria = "Limón"
$variedad = "Sin especificar"
$origen = "Provincia de Limarí"

# Row 327: 1a amarillo
$ws.Cells.Item(327, 1).Value = $mercadoId
$ws.Cells.Item(327, 2).Value = $mercado
$ws.Cells.Item(327, 3).Value = $region
$ws.Cells.Item(327, 4).Value = 44524
$ws.Cells.Item(327, 5).Value = $codreg
$ws.Cells.Item(327, 6).Value = $tipo
$ws.Cells.Item(327, 7).Value = $productoId
$ws.Cells.Item(327, 8).Value = $producto
$ws.Cells.Item(327, 9).Value = $categoriaId
$ws.Cells.Item(327, 10).Value = $categoria
$ws.Cells.Item(327, 11).Value = $variedad
$ws.Cells.Item(327, 12).Value = "1a amarillo"
$ws.Cells.Item(327, 13).Value = 900
$ws.Cells.Item(327, 14).Value = 5800
$ws.Cells.Item(327, 15).Value = 6000
$ws.Cells.Item(327, 16).Value = 5900
$ws.Cells.Item(327, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(327, 18).Value = $origen
$ws.Cells.Item(327, 19).Value = 369
$ws.Cells.Item(327, 20).Value = 16

# Row 328: 2a amarillo
$ws.Cells.Item(328, 1).Value = $mercadoId
$ws.Cells.Item(328, 2).Value = $mercado
$ws.Cells.Item(328, 3).Value = $region
$ws.Cells.Item(328, 4).Value = 44524
$ws.Cells.Item(328, 5).Value = $codreg
$ws.Cells.Item(328, 6).Value = $tipo
$ws.Cells.Item(328, 7).Value = $productoId
$ws.Cells.Item(328, 8).Value = $producto
$ws.Cells.Item(328, 9).Value = $categoriaId
$ws.Cells.Item(328, 10).Value = $categoria
$ws.Cells.Item(328, 11).Value = $variedad
$ws.Cells.Item(328, 12).Value = "2a amarillo"
$ws.Cells.Item(328, 13).Value = 750
$ws.Cells.Item(328, 14).Value = 4300
$ws.Cells.Item(328, 15).Value = 4500
$ws.Cells.Item(328, 16).Value = 4400
$ws.Cells.Item(328, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(328, 18).Value = $origen
$ws.Cells.Item(328, 19).Value = 275
$ws.Cells.Item(328, 20).Value = 16

# Row 329: 3a amarillo
$ws.Cells.Item(329, 1).Value = $mercadoId
$ws.Cells.Item(329, 2).Value = $mercado
$ws.Cells.Item(329, 3).Value = $region
$ws.Cells.Item(329, 4).Value = 44524
$ws.Cells.Item(329, 5).Value = $codreg
$ws.Cells.Item(329, 6).Value = $tipo
$ws.Cells.Item(329, 7).Value = $productoId
$ws.Cells.Item(329, 8).Value = $producto
$ws.Cells.Item(329, 9).Value = $categoriaId
$ws.Cells.Item(329, 10).Value = $categoria
$ws.Cells.Item(329, 11).Value = $variedad
$ws.Cells.Item(329, 12).Value = "3a amarillo"
$ws.Cells.Item(329, 13).Value = 510
$ws.Cells.Item(329, 14).Value = 2800
$ws.Cells.Item(329, 15).Value = 3000
$ws.Cells.Item(329, 16).Value = 2900
$ws.Cells.Item(329, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(329, 18).Value = $origen
$ws.Cells.Item(329, 19).Value = 181
$ws.Cells.Item(329, 20).Value = 16

# Make sure the date column keeps its date number format (style carried
# over from the Insert() above already covers D327:D329, this is just
# a safety net in case the insert didn't propagate it).
$ws.Range("D327:D329").NumberFormat = "YYYY-MM-DD HH:MM:SS"
